# Credentials.xlsx - "Made changes in the data file"
# Row 4 (Dokuparthi) data updated:
#   B4: Dokuparthi@gmail.com -> dokuparthi@gmail.com (lowercase)
#   C4: 12345 -> 123Dokuparthi (now text)
#   D4: 12345 -> 123fdsdfdf (now text)
# Selection moved to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

$ws.Range("B4").Value = "dokuparthi@gmail.com"
$ws.Range("C4").Value = "123Dokuparthi"
$ws.Range("D4").Value = "123fdsdfdf"

[void]$ws.Range("C3").Select()
